$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B21:G21").Copy()
$ws.Range("B22:G22").PasteSpecial(-4122)

$ws.Range("B22").Value = "EDENILSON SILVA"
$ws.Range("C22").Value = "e26d1cd1918b4c7d99e4509543ea983a"
$ws.Range("D22").Value = 44852
$ws.Range("E22").Value = 365
$ws.Range("F22").Value = "-"
$ws.Range("G22").Value = "VENDA 18 (18/10)"
